# Add two new columns, I (I0) and J (IF), to the stats table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - headers inherit the same bold/bordered style as the
# existing headers (s="1"), so copy H1's formatting onto I1:J1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4

$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 6

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 4

$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 6

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 2
